$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = "L682801"
$ws.Range("C5").Value = "SB#5"
$ws.Range("E5").Value = 1260
$ws.Range("F5").Value = "T"
$ws.Range("H5").Value = 45130.041998113426
$ws.Range("J5").Value = "07/11/23 17:56"
$ws.Range("K5").Value = "07/11/23 17:56"
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = "`$1,280 as of 7/11/2023 11:50:37 AM"
$ws.Range("N5").Value = 1280

# Row 6
$ws.Range("A6").Value = "LK644532"
$ws.Range("C6").Value = "SCL ENTERPRISES LAUNDRY"
$ws.Range("E6").Value = 1320
$ws.Range("F6").Value = "T"
$ws.Range("H6").Value = 45155.041998113426
$ws.Range("J6").Value = "07/11/23 21:44"
$ws.Range("K6").Value = "07/11/23 21:44"
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = "`$1,360 as of 7/8/2023 7:29:51 PM"
$ws.Range("N6").Value = 1360

# Row 7
$ws.Range("A7").Value = "L474761"
$ws.Range("C7").Value = "BABS MARKET"
$ws.Range("E7").Value = 1520
$ws.Range("F7").Value = "T"
$ws.Range("H7").Value = 45126.041998113426
$ws.Range("J7").Value = "07/11/23 12:57"
$ws.Range("K7").Value = "07/11/23 12:57"
$ws.Range("L7").Value = 40
$ws.Range("M7").Value = "`$1,520 as of 7/11/2023 10:57:39 AM"
$ws.Range("N7").Value = 1600

# Row 9
$ws.Range("E9").Value = 2420
$ws.Range("H9").Value = 45131.041998113426
$ws.Range("J9").Value = "07/11/23 16:50"
$ws.Range("K9").Value = "07/11/23 16:50"
$ws.Range("M9").Value = "`$2,440 as of 7/11/2023 11:25:19 AM"
$ws.Range("N9").Value = 2440

# Row 10
$ws.Range("A10").Value = "L474792"
$ws.Range("C10").Value = "NICK SHELL SERVICE"
$ws.Range("E10").Value = 2580
$ws.Range("F10").Value = "T"
$ws.Range("H10").Value = 45164.041998113426
$ws.Range("I10").ClearContents()
$ws.Range("J10").Value = "07/10/23 22:05"
$ws.Range("K10").Value = "07/10/23 22:05"
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = "`$2,580 as of 7/10/2023 8:05:11 PM"
$ws.Range("N10").Value = 2600

# Row 11
$ws.Range("A11").Value = "LK236828"
$ws.Range("C11").Value = "WORLDWIDE AUTOMOTIVE"
$ws.Range("E11").Value = 2640
$ws.Range("F11").Value = "T"
$ws.Range("H11").Value = 45135.041998113426
$ws.Range("J11").Value = "07/11/23 13:12"
$ws.Range("K11").Value = "07/11/23 13:12"
$ws.Range("L11").Value = 80
$ws.Range("M11").Value = "`$2,640 as of 7/11/2023 11:12:32 AM"
$ws.Range("N11").Value = 2660

# Row 12
$ws.Range("A12").Value = "L474746"
$ws.Range("C12").Value = "ZACATES MARKET"
$ws.Range("E12").Value = 2680
$ws.Range("F12").Value = "T"
$ws.Range("H12").Value = 45129.041998113426
$ws.Range("J12").Value = "07/11/23 21:00"
$ws.Range("K12").Value = "07/11/23 15:31"
$ws.Range("M12").Value = "`$2,780 as of 7/10/2023 5:55:27 PM"
$ws.Range("N12").Value = 2680

# Row 13
$ws.Range("A13").Value = "LK561655"
$ws.Range("C13").Value = "CRENSHAW CRAVOR #2"
$ws.Range("E13").Value = 2780
$ws.Range("F13").Value = "T"
$ws.Range("H13").ClearContents()
$ws.Range("I13").Value = "ATM Inactive greater than 48 minutes"
$ws.Range("J13").Value = "01/23/20 08:24"
$ws.Range("K13").Value = "01/23/20 08:24"
$ws.Range("M13").Value = "`$2,780 as of 1/23/2020 6:24:32 AM"
$ws.Range("N13").Value = 2800

# Row 14
$ws.Range("A14").Value = "L688961"
$ws.Range("C14").Value = "MONA MART"
$ws.Range("E14").Value = 2860
$ws.Range("F14").Value = "T"
$ws.Range("H14").Value = 45143.041998113426
$ws.Range("I14").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J14").Value = "07/09/23 19:27"
$ws.Range("K14").Value = "07/09/23 19:27"
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = "`$2,860 as of 7/9/2023 5:27:48 PM"
$ws.Range("N14").Value = 2920

# Row 16
$ws.Range("E16").Value = 3980
$ws.Range("H16").Value = 45127.041998113426
$ws.Range("J16").Value = "07/11/23 21:55"
$ws.Range("K16").Value = "07/11/23 19:13"
$ws.Range("M16").Value = "`$4,040 as of 7/11/2023 11:33:08 AM"
$ws.Range("N16").Value = 3980

# Row 17
$ws.Range("A17").Value = "L475090"
$ws.Range("C17").Value = "S.B. 2"
$ws.Range("E17").Value = 4960
$ws.Range("F17").Value = "T"
$ws.Range("H17").Value = 45123.041998113426
$ws.Range("J17").Value = "07/11/23 21:37"
$ws.Range("K17").Value = "07/11/23 21:25"
$ws.Range("L17").Value = 20
$ws.Range("M17").Value = "`$5,580 as of 7/11/2023 9:25:53 AM"
$ws.Range("N17").Value = 4960

# Row 18
$ws.Range("A18").Value = "L704741"
$ws.Range("C18").Value = "W ADAMS COIN LAUNDRY"
$ws.Range("E18").Value = 4980
$ws.Range("F18").Value = "T"
$ws.Range("H18").Value = 45125.041998113426
$ws.Range("J18").Value = "07/12/23 00:25"
$ws.Range("K18").Value = "07/12/23 00:25"
$ws.Range("M18").Value = "`$5,740 as of 7/11/2023 11:50:13 AM"
$ws.Range("N18").Value = 5000

# Row 19
$ws.Range("A19").Value = "L678988"
$ws.Range("C19").Value = "PAYELESS MARKET"
$ws.Range("E19").Value = 5560
$ws.Range("F19").Value = "T"
$ws.Range("H19").Value = 45142.041998113426
$ws.Range("J19").Value = "07/11/23 19:48"
$ws.Range("K19").Value = "07/11/23 19:48"
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = "`$5,720 as of 7/11/2023 11:38:49 AM"
$ws.Range("N19").Value = 5660

# Row 20
$ws.Range("A20").Value = "L488595"
$ws.Range("C20").Value = "N S MART"
$ws.Range("E20").Value = 5600
$ws.Range("F20").Value = "T"
$ws.Range("H20").Value = 45259.041998113426
$ws.Range("J20").Value = "07/11/23 22:35"
$ws.Range("K20").Value = "07/11/23 22:35"
$ws.Range("L20").Value = 60
$ws.Range("M20").Value = "`$5,720 as of 7/9/2023 9:11:13 PM"
$ws.Range("N20").Value = 5700

# Row 21
$ws.Range("A21").Value = "LK864765"
$ws.Range("C21").Value = "SKY LIQUOR"
$ws.Range("E21").Value = 5860
$ws.Range("F21").Value = "T"
$ws.Range("H21").Value = 45128.041998113426
$ws.Range("J21").Value = "07/11/23 23:52"
$ws.Range("K21").Value = "07/11/23 21:55"
$ws.Range("M21").Value = "`$6,200 as of 7/10/2023 7:01:11 PM"
$ws.Range("N21").Value = 5860

# Row 22
$ws.Range("A22").Value = "L688966"
$ws.Range("C22").Value = "LACON MINI MART"
$ws.Range("E22").Value = 6400
$ws.Range("F22").Value = "T"
$ws.Range("H22").Value = 45343.041998113426
$ws.Range("I22").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J22").Value = "07/09/23 22:56"
$ws.Range("K22").Value = "07/09/23 15:28"
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = "`$6,400 as of 7/9/2023 1:28:46 PM"
$ws.Range("N22").Value = 6420

# Row 23
$ws.Range("A23").Value = "LK923383"
$ws.Range("C23").Value = "SAMYS PHONE CARDS"
$ws.Range("E23").Value = 6480
$ws.Range("F23").Value = "T"
$ws.Range("H23").Value = 45131.041998113426
$ws.Range("J23").Value = "07/11/23 22:42"
$ws.Range("K23").Value = "07/11/23 22:42"
$ws.Range("L23").Value = 100
$ws.Range("M23").Value = "`$7,420 as of 7/11/2023 11:13:23 AM"
$ws.Range("N23").Value = 6680

# Row 24
$ws.Range("A24").Value = "L474817"
$ws.Range("C24").Value = "SAFETY MARKET"
$ws.Range("E24").Value = 7240
$ws.Range("F24").Value = "T"
$ws.Range("H24").Value = 45134.041998113426
$ws.Range("J24").Value = "07/11/23 13:51"
$ws.Range("K24").Value = "07/11/23 00:34"
$ws.Range("L24").Value = 120
$ws.Range("M24").Value = "`$7,240 as of 7/11/2023 11:51:07 AM"
$ws.Range("N24").Value = 7240

# Row 25
$ws.Range("A25").Value = "L697589"
$ws.Range("C25").Value = "S B DISCOUNT MART"
$ws.Range("E25").Value = 8380
$ws.Range("F25").Value = "T"
$ws.Range("H25").Value = 45124.041998113426
$ws.Range("I25").ClearContents()
$ws.Range("J25").Value = "07/11/23 22:38"
$ws.Range("K25").Value = "07/11/23 22:38"
$ws.Range("L25").Value = 40
$ws.Range("M25").Value = "`$9,200 as of 7/11/2023 11:53:30 AM"
$ws.Range("N25").Value = 8460

# Row 26
$ws.Range("A26").Value = "L697590"
$ws.Range("C26").Value = "S B MARKET ST"
$ws.Range("E26").Value = 8780
$ws.Range("F26").Value = "T"
$ws.Range("H26").ClearContents()
$ws.Range("I26").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J26").Value = "06/29/23 11:36"
$ws.Range("K26").Value = "06/29/23 11:36"
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = "`$8,780 as of 6/29/2023 9:36:36 AM"
$ws.Range("N26").Value = 8800

# Row 27
$ws.Range("E27").Value = 11780
$ws.Range("H27").Value = 45125.041998113426
$ws.Range("J27").Value = "07/11/23 22:44"
$ws.Range("K27").Value = "07/11/23 22:44"
$ws.Range("M27").Value = "`$12,320 as of 7/11/2023 11:23:11 AM"
$ws.Range("N27").Value = 11820

# Row 28
$ws.Range("E28").Value = 105800

